# Updated symbol list — refresh Price (column D) values for the rows
# whose quoted prices changed since the last scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-PriceText($row, $value) {
    # Column D = "Price". These cells store the price as TEXT (not a
    # number), so assign via a leading apostrophe to force text storage,
    # then restore the cell's original style so the quote-prefix flag
    # doesn't stick around as a visible formatting change.
    $cell = $ws.Cells.Item($row, 4)
    $origStyle = $cell.Style
    $cell.Value2 = "'" + $value
    $cell.Style = $origStyle
}

Set-PriceText 2 "286.46"
Set-PriceText 3 "21.25"
Set-PriceText 4 "6.457"
Set-PriceText 5 "0.06383"
Set-PriceText 6 "3.605"
Set-PriceText 7 "1.581"
Set-PriceText 8 "6.573"
Set-PriceText 9 "0.8262"
Set-PriceText 10 "0.01422"
Set-PriceText 11 "0.1679"
Set-PriceText 12 "0.08806"
Set-PriceText 13 "0.03700"
Set-PriceText 14 "0.03212"
Set-PriceText 15 "0.09193"
Set-PriceText 16 "3.711"
Set-PriceText 17 "0.001647"
Set-PriceText 18 "0.04764"
Set-PriceText 19 "0.006177"
Set-PriceText 20 "0.006292"
Set-PriceText 22 "0.0001603"
Set-PriceText 23 "3.780"
Set-PriceText 25 "0.3358"
Set-PriceText 26 "0.1262"
Set-PriceText 28 "0.0002711"
Set-PriceText 40 "0.04783"
Set-PriceText 41 "0.003541"
Set-PriceText 42 "0.003456"
Set-PriceText 43 "0.1121"
Set-PriceText 44 "0.01175"
Set-PriceText 45 "0.00006953"
Set-PriceText 47 "0.9351"
Set-PriceText 48 "0.008054"
Set-PriceText 49 "0.00001503"
